$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "FC Barcelona vs Rayo Vallecano / January 19th 2022 / 55" row (row 1) is
# removed entirely; every following row shifts up by one.
$ws.Rows.Item(1).Delete()

# Re-create column C (the score/attendance figures) for the remaining 7 rows so
# that the new values are written as plain text (matching how the other score
# cells in this sheet are stored), rather than being auto-typed as numbers.
$scores = $ws.Range("C1:C7")
$scores.Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

$scores.NumberFormat = "@"
$ws.Range("C1").Value = "81"
$ws.Range("C2").Value = "79"
$ws.Range("C3").Value = "90"
$ws.Range("C4").Value = "56"
$ws.Range("C5").Value = "62"
$ws.Range("C6").Value = "53"
$ws.Range("C7").Value = "59"
$scores.ClearFormats()
